# Auto-generated Excel COM-interop script
# Applies figure1_data.xlsx edits per commit 'Commands as sent to SSC on 23 November'

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows whose number_pass (C) and share_pass (D) values change
$changedRows = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,34,35,36,37,38,39,40,41,42,43,44,45,46,47,48,49,50,51,52,53,54,55,56,57,58,59,60,61,62,63,95,96,97,98,99,100,101,102,103,104,105,106,107,108,109,110,111,112,113,114,115,116,117,118,119,120,121,122,123,124,125)
$newC = @(136,135,135,135,134,134,134,133,133,133,133,131,130,127,127,127,127,126,123,118,118,117,110,109,105,104,102,102,96,95,14,137,137,137,137,136,136,136,133,132,131,129,99,98,97,95,94,91,89,87,84,82,76,13,9,9,5,5,3,1,1,1,75,75,75,75,75,75,70,70,70,70,69,65,65,65,65,65,65,64,64,64,63,63,63,63,63,63,62,62,62,62,46)
$newD = @(0.9577464788732394,0.9507042253521126,0.9507042253521126,0.9507042253521126,0.9436619718309859,0.9436619718309859,0.9436619718309859,0.9366197183098591,0.9366197183098591,0.9366197183098591,0.9366197183098591,0.9225352112676056,0.9154929577464789,0.8943661971830986,0.8943661971830986,0.8943661971830986,0.8943661971830986,0.8873239436619719,0.8661971830985915,0.8309859154929577,0.8309859154929577,0.823943661971831,0.7746478873239436,0.7676056338028169,0.7394366197183099,0.7323943661971831,0.7183098591549296,0.7183098591549296,0.676056338028169,0.6690140845070423,0.0985915492957746,0.9647887323943662,0.9647887323943662,0.9647887323943662,0.9647887323943662,0.9577464788732394,0.9577464788732394,0.9577464788732394,0.9366197183098591,0.9295774647887324,0.9225352112676056,0.9084507042253521,0.6971830985915493,0.6901408450704225,0.6830985915492958,0.6690140845070423,0.6619718309859155,0.6408450704225352,0.6267605633802817,0.6126760563380281,0.5915492957746479,0.5774647887323944,0.5352112676056338,0.0915492957746479,0.0633802816901408,0.0633802816901408,0.0352112676056338,0.0352112676056338,0.0211267605633803,0.0070422535211268,0.0070422535211268,0.0070422535211268,0.528169014084507,0.528169014084507,0.528169014084507,0.528169014084507,0.528169014084507,0.528169014084507,0.4929577464788732,0.4929577464788732,0.4929577464788732,0.4929577464788732,0.4859154929577465,0.4577464788732394,0.4577464788732394,0.4577464788732394,0.4577464788732394,0.4577464788732394,0.4577464788732394,0.4507042253521127,0.4507042253521127,0.4507042253521127,0.4436619718309859,0.4436619718309859,0.4436619718309859,0.4436619718309859,0.4436619718309859,0.4436619718309859,0.4366197183098591,0.4366197183098591,0.4366197183098591,0.4366197183098591,0.323943661971831)

for ($i = 0; $i -lt $changedRows.Length; $i++) {
    $r = $changedRows[$i]
    $ws.Cells.Item($r, 3).Value = $newC[$i]
    $ws.Cells.Item($r, 4).Value = $newD[$i]
}

# Rows for axiom groups 5 and 6 (rows 126-187) are removed entirely (all cells cleared)
$ws.Range("A126:D187").ClearContents()

Write-Host "Edit complete"
